$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of names
$ws.Range("A3").Value = "محمود"
$ws.Range("A4").Value = "علي"
$ws.Range("A5").Value = "خالد"
$ws.Range("A6").Value = "الاسم"
$ws.Range("B6").Value = "الاسم"

# Update the view: scroll back to A1 and select B6
$ws.Range("A1").Select()
$ws.Range("B6").Select()
